$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 197.5
$ws.Range("I4").Value = 200
$ws.Range("J4").Value = 195
$ws.Range("K4").Value = 200
$ws.Range("L4").Value = 195
$ws.Range("M4").Value = -86
$ws.Range("N4").Value = -423

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1110.4706
$ws.Range("I6").Value = 151.46153
$ws.Range("J6").Value = 4227.25
$ws.Range("K6").Value = 454.38459
$ws.Range("L6").Value = 12681.75
$ws.Range("M6").Value = -342.38459
$ws.Range("N6").Value = -12905.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1047.375
$ws.Range("I19").Value = 1163
$ws.Range("J19").Value = 978
$ws.Range("K19").Value = 1163
$ws.Range("L19").Value = 978
$ws.Range("M19").Value = -988
$ws.Range("N19").Value = -1328

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 6480.722
$ws.Range("I33").Value = 7515.857
$ws.Range("J33").Value = 2857.75
$ws.Range("K33").Value = 7515.857
$ws.Range("L33").Value = 2857.75
$ws.Range("M33").Value = -7286.857
$ws.Range("N33").Value = -3315.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 589.5714
$ws.Range("I92").Value = 581.36365
$ws.Range("J92").Value = 619.6667
$ws.Range("K92").Value = 581.36365
$ws.Range("L92").Value = 619.6667
$ws.Range("M92").Value = 666.63635
$ws.Range("N92").Value = -3115.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 35050.484
$ws.Range("I100").Value = 46360.824
$ws.Range("J100").Value = 2533.25
$ws.Range("K100").Value = 46360.824
$ws.Range("L100").Value = 2533.25
$ws.Range("M100").Value = -45819.824
$ws.Range("N100").Value = -3615.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 2132.7932
$ws.Range("I107").Value = 1660.64
$ws.Range("J107").Value = 5083.75
$ws.Range("K107").Value = 1660.64
$ws.Range("L107").Value = 5083.75
$ws.Range("M107").Value = 259.3599999999999
$ws.Range("N107").Value = -8923.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1648.3334
$ws.Range("I112").Value = 921.25
$ws.Range("J112").Value = 1748.6207
$ws.Range("K112").Value = 2763.75
$ws.Range("L112").Value = 5245.8621
$ws.Range("M112").Value = -1655.75
$ws.Range("N112").Value = -7461.8621

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1911.3334
$ws.Range("I132").Value = 1798.7778
$ws.Range("J132").Value = 2249
$ws.Range("K132").Value = 5396.3334
$ws.Range("L132").Value = 6747
$ws.Range("M132").Value = -2866.3334
$ws.Range("N132").Value = -11807

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3217.5
$ws.Range("I2").Value = 2460.9412
$ws.Range("J2").Value = 6432.875
$ws.Range("K2").Value = 2460.9412
$ws.Range("L2").Value = 6432.875
$ws.Range("M2").Value = -2347.9412
$ws.Range("N2").Value = -6658.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 11687.308
$ws.Range("I45").Value = 13035.7
$ws.Range("J45").Value = 7192.6665
$ws.Range("K45").Value = 13035.7
$ws.Range("L45").Value = 7192.6665
$ws.Range("M45").Value = -12658.7
$ws.Range("N45").Value = -7946.6665

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4928.64
$ws.Range("I61").Value = 2726.4736
$ws.Range("J61").Value = 11902.167
$ws.Range("K61").Value = 2726.4736
$ws.Range("L61").Value = 11902.167
$ws.Range("M61").Value = -2514.4736
$ws.Range("N61").Value = -12326.167

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 974.36365
$ws.Range("I97").Value = 1021.95
$ws.Range("J97").Value = 498.5
$ws.Range("K97").Value = 1021.95
$ws.Range("L97").Value = 498.5
$ws.Range("M97").Value = -525.95
$ws.Range("N97").Value = -1490.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 5385.9614
$ws.Range("I102").Value = 5385.9614
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 5385.9614
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -3763.9614

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 3217.5
$ws.Range("I116").Value = 2460.9412
$ws.Range("J116").Value = 6432.875
$ws.Range("K116").Value = 2460.9412
$ws.Range("L116").Value = 6432.875
$ws.Range("M116").Value = -166.9412000000002
$ws.Range("N116").Value = -11020.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3180.5557
$ws.Range("I122").Value = 3171.375
$ws.Range("J122").Value = 3254
$ws.Range("K122").Value = 9514.125
$ws.Range("L122").Value = 9762
$ws.Range("M122").Value = -7064.125
$ws.Range("N122").Value = -14662

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 5875.362
$ws.Range("I132").Value = 6168.886
$ws.Range("J132").Value = 1570.3334
$ws.Range("K132").Value = 18506.658
$ws.Range("L132").Value = 4711.0002
$ws.Range("M132").Value = -15976.658
$ws.Range("N132").Value = -9771.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4928.64
$ws.Range("I136").Value = 2726.4736
$ws.Range("J136").Value = 11902.167
$ws.Range("K136").Value = 8179.4208
$ws.Range("L136").Value = 35706.501
$ws.Range("M136").Value = -5629.4208
$ws.Range("N136").Value = -40806.501

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3217.5
$ws.Range("I3").Value = 2460.9412
$ws.Range("J3").Value = 6432.875
$ws.Range("K3").Value = 2460.9412
$ws.Range("L3").Value = 6432.875
$ws.Range("M3").Value = -2346.9412
$ws.Range("N3").Value = -6660.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 9517.191999999999
$ws.Range("I134").Value = 6611.4546
$ws.Range("J134").Value = 25498.75
$ws.Range("K134").Value = 19834.3638
$ws.Range("L134").Value = 76496.25
$ws.Range("M134").Value = -17299.3638
$ws.Range("N134").Value = -81566.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 2413.7144
$ws.Range("I15").Value = 744
$ws.Range("J15").Value = 3081.6
$ws.Range("K15").Value = 744
$ws.Range("L15").Value = 3081.6
$ws.Range("M15").Value = -574
$ws.Range("N15").Value = -3421.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3233
$ws.Range("I31").Value = 2139.0476
$ws.Range("J31").Value = 5785.5557
$ws.Range("K31").Value = 2139.0476
$ws.Range("L31").Value = 5785.5557
$ws.Range("M31").Value = -1844.0476
$ws.Range("N31").Value = -6375.5557

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3233
$ws.Range("I34").Value = 2139.0476
$ws.Range("J34").Value = 5785.5557
$ws.Range("K34").Value = 2139.0476
$ws.Range("L34").Value = 5785.5557
$ws.Range("M34").Value = -1937.0476
$ws.Range("N34").Value = -6189.5557

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3613.08
$ws.Range("I58").Value = 2131.8
$ws.Range("J58").Value = 9538.200000000001
$ws.Range("K58").Value = 2131.8
$ws.Range("L58").Value = 9538.200000000001
$ws.Range("M58").Value = -1928.8
$ws.Range("N58").Value = -9944.200000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2323.125
$ws.Range("I122").Value = 2264.1667
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 6792.500100000001
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -4342.500100000001
$ws.Range("N122").Value = -12400

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 27811.445
$ws.Range("I132").Value = 17782.473
$ws.Range("J132").Value = 47869.39
$ws.Range("K132").Value = 53347.41900000001
$ws.Range("L132").Value = 143608.17
$ws.Range("M132").Value = -50817.41900000001
$ws.Range("N132").Value = -148668.17

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3613.08
$ws.Range("I136").Value = 2131.8
$ws.Range("J136").Value = 9538.200000000001
$ws.Range("K136").Value = 6395.400000000001
$ws.Range("L136").Value = 28614.6
$ws.Range("M136").Value = -3845.400000000001
$ws.Range("N136").Value = -33714.60000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2106.0286
$ws.Range("I5").Value = 921.7
$ws.Range("J5").Value = 2579.76
$ws.Range("K5").Value = 2765.1
$ws.Range("L5").Value = 7739.280000000001
$ws.Range("M5").Value = -2653.1
$ws.Range("N5").Value = -7963.280000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H93").Value = 7250
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 7250
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 21750
$ws.Range("N93").Value = -25494

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 10001494
$ws.Range("I122").Value = 1992.4
$ws.Range("J122").Value = 20000996
$ws.Range("K122").Value = 17931.6
$ws.Range("L122").Value = 180008964
$ws.Range("M122").Value = -15481.6
$ws.Range("N122").Value = -180013864

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 2106.0286
$ws.Range("I135").Value = 921.7
$ws.Range("J135").Value = 2579.76
$ws.Range("K135").Value = 8295.300000000001
$ws.Range("L135").Value = 23217.84
$ws.Range("M135").Value = -5760.300000000001
$ws.Range("N135").Value = -28287.84

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 18502
$ws.Range("I5").Value = 18502
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 18502
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -18390

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5201.75
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 5415.2
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 16245.6
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -21185.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1381.409
$ws.Range("I46").Value = 1016.5
$ws.Range("J46").Value = 1518.25
$ws.Range("K46").Value = 1016.5
$ws.Range("L46").Value = 1518.25
$ws.Range("M46").Value = -828.5
$ws.Range("N46").Value = -1894.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4166.6665
$ws.Range("I122").Value = 3995
$ws.Range("J122").Value = 4252.5
$ws.Range("K122").Value = 11985
$ws.Range("L122").Value = 12757.5
$ws.Range("M122").Value = -9535
$ws.Range("N122").Value = -17657.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 49999
$ws.Range("I42").Value = 49999
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 49999
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -49621
$ws.Range("N42").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 24009.098
$ws.Range("I132").Value = 16622.178
$ws.Range("J132").Value = 39919.383
$ws.Range("K132").Value = 49866.534
$ws.Range("L132").Value = 119758.149
$ws.Range("M132").Value = -47336.534
$ws.Range("N132").Value = -124818.149

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 53365.58
$ws.Range("I136").Value = 53365.58
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 160096.74
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -157546.74
